# Fruta / hortaliza, semanal
# Insert 3 new weekly rows for "Vega Monumental Concepción - Plátano" right
# before the existing row 490, pushing the rest of the historical rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 490 downward by 3 (Excel copies formatting from the row above
# into the freshly inserted blank rows, same as a manual "Insert" in the UI).
$ws.Rows("490:492").Insert()

# Row 490 - Maduro
$ws.Cells.Item(490, 1).Value  = 11
$ws.Cells.Item(490, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(490, 3).Value  = "Bíobío"
$ws.Cells.Item(490, 4).Value  = "2022-07-05"
$ws.Cells.Item(490, 5).Value  = 8
$ws.Cells.Item(490, 6).Value  = "Fruta"
$ws.Cells.Item(490, 7).Value  = 100108
$ws.Cells.Item(490, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(490, 9).Value  = 100108006
$ws.Cells.Item(490, 10).Value = "Plátano"
$ws.Cells.Item(490, 11).Value = "Sin especificar"
$ws.Cells.Item(490, 12).Value = "Maduro"
$ws.Cells.Item(490, 13).Value = 100
$ws.Cells.Item(490, 14).Value = 20000
$ws.Cells.Item(490, 15).Value = 20000
$ws.Cells.Item(490, 16).Value = 20000
$ws.Cells.Item(490, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(490, 18).Value = "Ecuador"
$ws.Cells.Item(490, 19).Value = 1000
$ws.Cells.Item(490, 20).Value = 20

# Row 491 - Pintón
$ws.Cells.Item(491, 1).Value  = 11
$ws.Cells.Item(491, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(491, 3).Value  = "Bíobío"
$ws.Cells.Item(491, 4).Value  = "2022-07-05"
$ws.Cells.Item(491, 5).Value  = 8
$ws.Cells.Item(491, 6).Value  = "Fruta"
$ws.Cells.Item(491, 7).Value  = 100108
$ws.Cells.Item(491, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(491, 9).Value  = 100108006
$ws.Cells.Item(491, 10).Value = "Plátano"
$ws.Cells.Item(491, 11).Value = "Sin especificar"
$ws.Cells.Item(491, 12).Value = "Pintón"
$ws.Cells.Item(491, 13).Value = 300
$ws.Cells.Item(491, 14).Value = 21000
$ws.Cells.Item(491, 15).Value = 21000
$ws.Cells.Item(491, 16).Value = 21000
$ws.Cells.Item(491, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(491, 18).Value = "Ecuador"
$ws.Cells.Item(491, 19).Value = 1050
$ws.Cells.Item(491, 20).Value = 20

# Row 492 - Primera Pintón
$ws.Cells.Item(492, 1).Value  = 11
$ws.Cells.Item(492, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(492, 3).Value  = "Bíobío"
$ws.Cells.Item(492, 4).Value  = "2022-07-05"
$ws.Cells.Item(492, 5).Value  = 8
$ws.Cells.Item(492, 6).Value  = "Fruta"
$ws.Cells.Item(492, 7).Value  = 100108
$ws.Cells.Item(492, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(492, 9).Value  = 100108006
$ws.Cells.Item(492, 10).Value = "Plátano"
$ws.Cells.Item(492, 11).Value = "Sin especificar"
$ws.Cells.Item(492, 12).Value = "Primera Pintón"
$ws.Cells.Item(492, 13).Value = 300
$ws.Cells.Item(492, 14).Value = 23000
$ws.Cells.Item(492, 15).Value = 23000
$ws.Cells.Item(492, 16).Value = 23000
$ws.Cells.Item(492, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(492, 18).Value = "Ecuador"
$ws.Cells.Item(492, 19).Value = 1150
$ws.Cells.Item(492, 20).Value = 20
